# Update the Tgfb1-Acvrl1 LR-pair sheet with refreshed TPM-derived values.
#
# The sheet lists, for every (Sending cluster, Target cluster) combination, the
# ligand statistics (columns G-J) that depend only on the Sending cluster, the
# receptor statistics (columns M-P) that depend only on the Target cluster, and
# the derived edge statistics (columns Q-T) that are simply the products of the
# matching ligand/receptor values:
#   Q = G*M   R = H*N   S = I*O   T = J*P
#
# New ligand values per Sending cluster (columns G,H,I,J):
$ligandBySending = @{
    "ECs"           = @(54.69462833333333, 164.083885,          0.2790924419198448,  0.2790924419198448)
    "FAPs"          = @(19.32115333333334, 57.96346000000001,   0.09859081282432611, 0.09859081282432611)
    "MuSCs"         = @(11.023718,         33.071154,           0.05625116157486912, 0.05625116157486911)
    "Resolving-Mac" = @(110.9336623333333, 332.800987,          0.5660655836809599,  0.5660655836809599)
}

# New receptor values per Target cluster (columns M,N,O,P):
$receptorByTarget = @{
    "ECs"           = @(31.40242733333333, 94.20728199999999, 0.5334014788811394, 0.5334014788811395)
    "FAPs"          = @(19.28977566666667, 57.869327,          0.327656036225058,  0.327656036225058)
    "MuSCs"         = @(1.868202333333333, 5.604607,           0.03173327580290011,0.03173327580290011)
    "Resolving-Mac" = @(6.311623666666667, 18.934871,          0.1072092090909023, 0.1072092090909024)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 17 }

for ($r = 2; $r -le $lastRow; $r++) {
    $sending = $ws.Range("A$r").Value()
    $target  = $ws.Range("D$r").Value()
    if (-not $sending -or -not $target) { continue }

    $lig = $ligandBySending[$sending]
    $rec = $receptorByTarget[$target]
    if (-not $lig -or -not $rec) { continue }

    $g = $lig[0]; $h = $lig[1]; $i = $lig[2]; $j = $lig[3]
    $m = $rec[0]; $n = $rec[1]; $o = $rec[2]; $p = $rec[3]

    $ws.Range("G$r").Value = $g
    $ws.Range("H$r").Value = $h
    $ws.Range("I$r").Value = $i
    $ws.Range("J$r").Value = $j

    $ws.Range("M$r").Value = $m
    $ws.Range("N$r").Value = $n
    $ws.Range("O$r").Value = $o
    $ws.Range("P$r").Value = $p

    $ws.Range("Q$r").Value = $g * $m
    $ws.Range("R$r").Value = $h * $n
    $ws.Range("S$r").Value = $i * $o
    $ws.Range("T$r").Value = $j * $p
}
